$d = $word.ActiveDocument

# Change 1: update years of experience in professional summary
$d.Content.Find.Execute(
    "21 years of experience",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "15+ years of experience", 2
)

# Change 2: enhance FLEEM bullet under Research Director role
$d.Content.Find.Execute(
    "Engineered FLEEM web application using Twilio's API to make thousands of simultaneous phone calls for IVR polls",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Conceived, architected, and engineered FLEEM web application using Twilio API handling tens of thousands of calls using emulated predictive dialer for regulated political surveys",
    2
)

# Change 3: add new bullet after the "Developed innovative approaches..." bullet
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`a") -eq "• Developed innovative approaches to visualizing demographic and market data, enhancing clients' understanding of research findings") {
        $p.Range.InsertAfter("`r• Trained staff on building Python tooling for report generation and analysis")
        break
    }
}
